$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Remove the "Meta description" paragraph near the top of the
#    document (it is being relocated near the end of the document).
# ------------------------------------------------------------------
$metaFind = $d.Content
$metaFound = $metaFind.Find.Execute("Meta description", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $metaFound) {
    throw "Could not find the 'Meta description' paragraph"
}
$metaFind.Expand(4)  # wdParagraph
$metaFind.Delete()

# ------------------------------------------------------------------
# 2. Insert a new bold paragraph ("Play Banana Splash Free - Review
#    of Online Slot Game") right before the final paragraph (the one
#    that used to hold the image-generation prompt).
# ------------------------------------------------------------------
$count = $d.Paragraphs.Count
$last = $d.Paragraphs.Item($count)
$insertPoint = $d.Range($last.Range.Start, $last.Range.Start)

$newParaXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Banana Splash Free - Review of Online Slot Game</w:t></w:r></w:p><w:p/></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$insertPoint.InsertXML($newParaXml)

# The InsertXML call above also introduces a stray empty paragraph
# (the trailing placeholder needed to keep the remaining content in
# its own paragraph) - remove it now that the split has happened.
$d.Paragraphs.Item($count + 1).Range.Delete()

# ------------------------------------------------------------------
# 3. Replace the final paragraph's text (the old image-generation
#    prompt) with the meta-description copy text.
# ------------------------------------------------------------------
$oldPrompt = "Create a feature image for Banana Splash that features a happy Maya warrior with glasses in cartoon style. The background should be a beach with a blue sea and a bright sun shining over the warrior's head. The warrior should be holding a banana in one hand and a beach ball in the other, with a big smile on their face. The other symbols from the game, such as a pineapple with a lifebuoy, a melon playing bongos, and a watermelon sunbathing, should be in the background, scattered on the beach. This feature image should capture the fun and cheerful theme of the game and entice players to try it out."
$newCopy = "Learn more about Banana Splash - a fun and simple online slot game with colorful graphics. Play Banana Splash for free and enjoy its bonus rounds."

$d.Content.Find.Execute($oldPrompt, $true, $false, $false, $false, $false, $true, 1, $false, $newCopy, 2) | Out-Null

Write-Output "done"
